# NIT-9011349394.xlsx update:
#  - Remove the previous account-statement (EC) data block and add the new one
#    (MONICA TAMAYO CASTAÑO, 10 periods) while keeping the first two workers.
#  - Refresh the totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: the data block grows from 3 rows (16:18) to 13 rows (16:28).
#        Insert 10 blank rows right after the first two (still-valid) records,
#        pushing the signature block (old rows 23/24) down to rows 33/34.
$ws.Rows("19:28").Insert()

# --- 2. Re-seed formatting for the grown block:
#        * row 18 held the "last row" (bottom-border) style -> move it to the
#          new last row (28) first, before it gets overwritten.
#        * row 16 holds the regular (interior) row style -> stamp it across
#          the newly inserted rows 18-27 so every interior row matches.
$ws.Range("B18:J18").Copy($ws.Range("B28:J28"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J27"))

# --- 3. Write the new table contents (rows 16-28).
$rows = @(
    @{R=16; C="1007855319"; D="ARAMIS MENDOZA VALDEZ"; E="1901"; F=16562; G=828116},
    @{R=17; C="1047422179"; D="YEISER EDUARDO AVILA HEREDIA"; E="1901"; F=1104; G=828116},
    @{R=18; C="73213306";   D="HAILER DE LA ROSA MONTIEL";    E="1807"; F=9375; G=781242},
    @{R=19; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2102"; F=35112; G=877803},
    @{R=20; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2101"; F=35112; G=877803},
    @{R=21; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2012"; F=35112; G=877803},
    @{R=22; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2011"; F=35112; G=877803},
    @{R=23; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2010"; F=35112; G=877803},
    @{R=24; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2009"; F=35112; G=877803},
    @{R=25; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2008"; F=35112; G=877803},
    @{R=26; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2007"; F=35112; G=877803},
    @{R=27; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2006"; F=35112; G=877803},
    @{R=28; C="43519583";   D="MONICA TAMAYO CASTAÑO";        E="2005"; F=35112; G=877803}
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}

# --- 4. Refresh the summary totals.
$ws.Range("E11").Value = 378161
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 12
